$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8, shifting "registro de marca/logo" row (old 8)
# and the blank row (old 9) down to rows 9 and 10.
$ws.Rows("8:8").Insert()

# --- Column C "costo" values for the existing expense rows ---
$ws.Range("C3").Value = "costo (USD)"
$ws.Range("C4").Value = 800
$ws.Range("C5").Value = "6500 - 13000"
$ws.Range("C6").Value = "100-300"

# Row 7 "dominio" cost: a date formatted as mmm-yy (Oct-20)
$ws.Range("C7").Value = 44105
$ws.Range("C7").NumberFormat = "mmm-yy"

# --- New row 8: "internet" ---
$ws.Range("B8").Value = "internet"
$ws.Range("C8").Value = "100-300"

# Row 9 now holds "registro de marca/logo" (shifted from row 8) - add its cost
$ws.Range("C9").Value = 1000

# --- New "ingresos" (income) block ---
$ws.Range("B12").Value = "ingresos"

$ws.Range("B13").Value = "proyecto pequenio"
$ws.Range("B14").Value = "proyecto mediano"
$ws.Range("B15").Value = "proyecto grande"

$ws.Range("C13").Value = "5k "
$ws.Range("C14").Value = "30k"
$ws.Range("C15").Value = "75k"

$ws.Range("E14").Value = $null
$ws.Range("E14").Style = $ws.Range("B10").Style

# Restore selection to match the final authored state
$ws.Range("F9").Select()
